# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity, Temperature and
# Proximity sheets, matching the new entries recorded for 2026-02-01
# around 18:33-18:34.

$wb = $excel.ActiveWorkbook

function Add-LogRow {
    param(
        $ws,
        [int]$row,
        [string]$date,
        [string]$timestamp,
        [string]$hour,
        [string]$location,
        [string]$value,
        [string]$status
    )

    # The Date column (and percentage-looking Values, e.g. "80.0%")
    # get silently auto-parsed into numbers by the Value setter (just
    # like typing them into real Excel), so force those particular
    # cells to Text format first to keep them as literal strings.
    $dateCell = $ws.Cells.Item($row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $date

    $ws.Cells.Item($row, 2).Value = $timestamp
    $ws.Cells.Item($row, 3).Value = $hour
    $ws.Cells.Item($row, 4).Value = $location

    $valueCell = $ws.Cells.Item($row, 5)
    if ($value -like "*%") {
        $valueCell.NumberFormat = "@"
    }
    $valueCell.Value = $value

    $ws.Cells.Item($row, 6).Value = $status
}

# --- PIR sheet: rows 65-68 -------------------------------------------------
$wsPir = $wb.Worksheets.Item("PIR")
Add-LogRow $wsPir 65 "2026-02-01" "18:33:48" "18:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPir 66 "2026-02-01" "18:33:53" "18:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPir 67 "2026-02-01" "18:33:58" "18:00" "Bathroom" "No Motion" "Inactive"
Add-LogRow $wsPir 68 "2026-02-01" "18:34:03" "18:00" "Bathroom" "No Motion" "Inactive"

# --- Humidity sheet: rows 123-127 ------------------------------------------
$wsHumidity = $wb.Worksheets.Item("Humidity")
Add-LogRow $wsHumidity 123 "2026-02-01" "18:33:43" "18:00" "Bathroom" "80.0%" "Active"
Add-LogRow $wsHumidity 124 "2026-02-01" "18:33:49" "18:00" "Bathroom" "81.1%" "Active"
Add-LogRow $wsHumidity 125 "2026-02-01" "18:33:54" "18:00" "Bathroom" "80.1%" "Active"
Add-LogRow $wsHumidity 126 "2026-02-01" "18:33:59" "18:00" "Bathroom" "81.1%" "Active"
Add-LogRow $wsHumidity 127 "2026-02-01" "18:34:04" "18:00" "Bathroom" "80.1%" "Active"

# --- Temperature sheet: rows 123-127 ---------------------------------------
$wsTemperature = $wb.Worksheets.Item("Temperature")
Add-LogRow $wsTemperature 123 "2026-02-01" "18:33:44" "18:00" "Bathroom" "29.1C" "Active"
Add-LogRow $wsTemperature 124 "2026-02-01" "18:33:49" "18:00" "Bathroom" "29.1C" "Active"
Add-LogRow $wsTemperature 125 "2026-02-01" "18:33:54" "18:00" "Bathroom" "29.1C" "Active"
Add-LogRow $wsTemperature 126 "2026-02-01" "18:33:59" "18:00" "Bathroom" "29.1C" "Active"
Add-LogRow $wsTemperature 127 "2026-02-01" "18:34:05" "18:00" "Bathroom" "29.1C" "Active"

# --- Proximity sheet: row 49 ------------------------------------------------
$wsProximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $wsProximity 49 "2026-02-01" "18:33:58" "18:00" "Bathroom Door" "EXIT" "User EXITED Bathroom"
